$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.141.69"
$ws.Range("D3").Value = "1.850.60"
$ws.Range("E3").Value = "  -2.26%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D5").Value = "'0.7039"
$ws.Range("E5").Value = "  -5.43%  "
$ws.Range("D6").Value = "'238.27"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.3048"
$ws.Range("E8").Value = "  -3.98%  "
$ws.Range("D9").Value = "'0.07486"
$ws.Range("E9").Value = "  +3.35%  "
$ws.Range("E10").Value = "  -6.79%  "
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7250"
$ws.Range("E12").Value = "  -5.10%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.219"
$ws.Range("E13").Value = "  -4.34%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.826.83"
$ws.Range("E14").Value = "  -3.77%  "
$ws.Range("D15").Value = "'89.11"
$ws.Range("E15").Value = "  -4.22%  "
$ws.Range("D16").Value = "29.254.74"
$ws.Range("E16").Value = "  -3.07%  "
$ws.Range("D17").Value = "'5.788"
$ws.Range("E17").Value = "  -6.44%  "
$ws.Range("D18").Value = "'239.05"
$ws.Range("E18").Value = "  -4.72%  "
$ws.Range("D19").Value = "'0.000007672"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").Value = "'13.06"
$ws.Range("E20").Value = "  -4.55%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "2.117.51"
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'7.565"
$ws.Range("E24").Value = "  -5.65%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'161.87"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").Value = "'8.992"
$ws.Range("E26").Value = "  -3.51%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.1462"
$ws.Range("E27").Value = "  -7.92%  "
$ws.Range("E28").Value = "  -4.28%  "
$ws.Range("D29").Value = "'1.939"
$ws.Range("E29").Value = "  -6.63%  "
$ws.Range("D30").Value = "'1.385"
$ws.Range("E30").Value = "  -6.32%  "
$ws.Range("D31").Value = "'4.566"
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").Value = "'1.494"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").Value = "'4.001"
$ws.Range("E33").Value = "  -5.42%  "
$ws.Range("D34").Value = "'0.05155"
$ws.Range("E34").Value = "  -4.75%  "
$ws.Range("D35").Value = "'1.188"
$ws.Range("E35").Value = "  -5.37%  "
$ws.Range("D36").Value = "'1.040"
$ws.Range("E36").Value = "  +4.60%  "
$ws.Range("D37").Value = "'0.7013"
$ws.Range("E37").Value = "  -8.75%  "
$ws.Range("D38").Value = "'2.644"
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("D39").Value = "'0.01869"
$ws.Range("E39").Value = "  -5.38%  "
$ws.Range("D40").Value = "'2.679"
$ws.Range("E40").Value = "  -3.42%  "
$ws.Range("D41").Value = "'0.9498"
$ws.Range("E41").Value = "  +9.08%  "
$ws.Range("D42").Value = "'6.011"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").Value = "1.080.04"
$ws.Range("E43").Value = "  -2.00%  "
$ws.Range("D44").Value = "'0.4302"
$ws.Range("E44").Value = "  -6.06%  "
$ws.Range("D45").Value = "'69.90"
$ws.Range("E45").Value = "  -4.19%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "'102.23"
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.019.42"
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.751"
$ws.Range("E49").Value = "  -6.55%  "
$ws.Range("D50").Value = "'9.179"
$ws.Range("E50").Value = "  -4.74%  "
$ws.Range("D51").Value = "'7.058"
$ws.Range("E51").Value = "  -7.55%  "
